$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change), and a row swap for
# RocketPoolETH / Quant (rows 49-50), per the commit diff.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.861.65"
$ws.Range("E2").Value = "  -0.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.779.63"
$ws.Range("E3").Value = "  -1.45%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.37"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.547"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.73"
$ws.Range("E8").Value = "  -2.44%  "

$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0676"
$ws.Range("E10").Value = "  -6.24%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.028.66"
$ws.Range("E12").Value = "  -1.75%  "

$ws.Range("E13").Value = "  +1.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.773.21"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.867.42"
$ws.Range("E15").Value = "  -1.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.609"
$ws.Range("E16").Value = "  -3.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.13"
$ws.Range("E17").Value = "  -2.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.77"
$ws.Range("E18").Value = "  -2.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "238.77"
$ws.Range("E19").Value = "  -3.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("E22").Value = "  -4.35%  "

$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -2.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.82"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.03"
$ws.Range("E26").Value = "  -1.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.07"
$ws.Range("E27").Value = "  -3.44%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("E31").Value = "  -3.60%  "

$ws.Range("E32").Value = "  -4.17%  "

$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("E34").Value = "  -2.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.388.40"
$ws.Range("E35").Value = "  -2.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.633"
$ws.Range("E36").Value = "  -3.19%  "

$ws.Range("E37").Value = "  -1.87%  "

$ws.Range("E38").Value = "  -1.36%  "

$ws.Range("E39").Value = "  +4.05%  "

$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.911"
$ws.Range("E41").Value = "  -3.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.21"
$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("E43").Value = "  -3.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.45"
$ws.Range("E44").Value = "  +11.62%  "

$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("E46").Value = "  +10.38%  "

$ws.Range("E47").Value = "  +2.49%  "

$ws.Range("E48").Value = "  -1.90%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.936.57"
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.58"
$ws.Range("E50").Value = "  -2.55%  "

$ws.Range("E51").Value = "  +0.12%  "
